# Applies the data shift for "Fruta, Feria Lagunitas de Puerto Montt - Kiwi"
# Rows 535-575 shift down by 3 (to 538-578); rows 535-537 get new weekly data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
  @{ Row = 535; D = 45106; L = 'Especial'; M = 200; N = 18000; O = 18000; P = 18000; Q = '$/caja 15 kilos'; R = 'Región de O''Higgins'; S = 1200 },
  @{ Row = 536; D = 45106; L = 'Primera'; M = 200; N = 15000; O = 15000; P = 15000; Q = '$/caja 15 kilos'; R = 'Región de O''Higgins'; S = 1000 },
  @{ Row = 537; D = 45106; L = 'Segunda'; M = 200; N = 13000; O = 13000; P = 13000; Q = '$/caja 15 kilos'; R = 'Región de O''Higgins'; S = 867 },
  @{ Row = 538; D = 44383; L = 'Especial'; M = 150; N = 15000; O = 15000; P = 15000; Q = '$/caja 15 kilos'; R = 'Región de O''Higgins'; S = 1000 },
  @{ Row = 539; D = 44383; L = 'Primera'; M = 300; N = 13000; O = 13500; P = 13250; Q = '$/caja 15 kilos'; R = 'Región de O''Higgins'; S = 883 },
  @{ Row = 540; D = 45068; L = 'Especial'; M = 200; N = 19000; O = 19000; P = 19000; Q = '$/caja 15 kilos'; R = 'Región de O''Higgins'; S = 1267 },
  @{ Row = 541; D = 45068; L = 'Primera'; M = 200; N = 17000; O = 17000; P = 17000; Q = '$/caja 15 kilos'; R = 'Región de O''Higgins'; S = 1133 },
  @{ Row = 542; D = 45068; L = 'Segunda'; M = 200; N = 14000; O = 14000; P = 14000; Q = '$/caja 15 kilos'; R = 'Región de O''Higgins'; S = 933 },
  @{ Row = 543; D = 44994; L = 'Especial'; M = 200; N = 25000; O = 25000; P = 25000; Q = '$/caja 15 kilos'; R = 'Región de O''Higgins'; S = 1667 },
  @{ Row = 544; D = 44994; L = 'Primera'; M = 400; N = 18000; O = 19000; P = 18500; Q = '$/caja 15 kilos'; R = 'Región de O''Higgins'; S = 1233 },
  @{ Row = 545; D = 44533; L = 'Especial'; M = 150; N = 22000; O = 22000; P = 22000; Q = '$/caja 15 kilos'; R = 'Provincia de Curicó'; S = 1467 },
  @{ Row = 546; D = 44533; L = 'Primera'; M = 300; N = 16000; O = 17000; P = 16500; Q = '$/caja 15 kilos'; R = 'Provincia de Curicó'; S = 1100 },
  @{ Row = 547; D = 44284; L = 'Primera'; M = 200; N = 16000; O = 17000; P = 16500; Q = '$/caja 15 kilos'; R = 'Región de O''Higgins'; S = 1100 },
  @{ Row = 548; D = 44756; L = 'Primera'; M = 200; N = 14000; O = 14000; P = 14000; Q = '$/caja 15 kilos'; R = 'Región de O''Higgins'; S = 933 },
  @{ Row = 549; D = 44756; L = 'Segunda'; M = 200; N = 12000; O = 12000; P = 12000; Q = '$/caja 15 kilos'; R = 'Región de O''Higgins'; S = 800 },
  @{ Row = 550; D = 44810; L = 'Especial'; M = 300; N = 17000; O = 17000; P = 17000; Q = '$/caja 15 kilos'; R = 'Región de O''Higgins'; S = 1133 },
  @{ Row = 551; D = 44810; L = 'Primera'; M = 300; N = 14000; O = 14000; P = 14000; Q = '$/caja 15 kilos'; R = 'Región de O''Higgins'; S = 933 },
  @{ Row = 552; D = 44810; L = 'Segunda'; M = 300; N = 12500; O = 12500; P = 12500; Q = '$/caja 15 kilos'; R = 'Región de O''Higgins'; S = 833 },
  @{ Row = 553; D = 44442; L = 'Especial'; M = 300; N = 21000; O = 21000; P = 21000; Q = '$/caja 15 kilos'; R = 'Provincia de Curicó'; S = 1400 },
  @{ Row = 554; D = 44442; L = 'Primera'; M = 150; N = 14000; O = 14000; P = 14000; Q = '$/caja 15 kilos'; R = 'Provincia de Curicó'; S = 933 },
  @{ Row = 555; D = 44523; L = 'Especial'; M = 350; N = 22000; O = 22000; P = 22000; Q = '$/caja 15 kilos'; R = 'Provincia de Curicó'; S = 1467 },
  @{ Row = 556; D = 44523; L = 'Primera'; M = 350; N = 18000; O = 18000; P = 18000; Q = '$/caja 15 kilos'; R = 'Provincia de Curicó'; S = 1200 },
  @{ Row = 557; D = 44523; L = 'Segunda'; M = 300; N = 14000; O = 14000; P = 14000; Q = '$/caja 15 kilos'; R = 'Provincia de Curicó'; S = 933 },
  @{ Row = 558; D = 44704; L = 'Segunda'; M = 300; N = 13000; O = 13000; P = 13000; Q = '$/caja 15 kilos'; R = 'Región de O''Higgins'; S = 867 },
  @{ Row = 559; D = 44504; L = 'Primera'; M = 200; N = 16000; O = 17000; P = 16500; Q = '$/caja 15 kilos'; R = 'Provincia de Curicó'; S = 1100 },
  @{ Row = 560; D = 45104; L = 'Especial'; M = 200; N = 18000; O = 18000; P = 18000; Q = '$/caja 15 kilos'; R = 'Región de O''Higgins'; S = 1200 },
  @{ Row = 561; D = 45104; L = 'Primera'; M = 200; N = 15000; O = 15000; P = 15000; Q = '$/caja 15 kilos'; R = 'Región de O''Higgins'; S = 1000 },
  @{ Row = 562; D = 45104; L = 'Segunda'; M = 200; N = 13000; O = 13000; P = 13000; Q = '$/caja 15 kilos'; R = 'Región de O''Higgins'; S = 867 },
  @{ Row = 563; D = 44859; L = 'Especial'; M = 200; N = 17000; O = 17000; P = 17000; Q = '$/caja 15 kilos'; R = 'Región de O''Higgins'; S = 1133 },
  @{ Row = 564; D = 44859; L = 'Primera'; M = 200; N = 15000; O = 15000; P = 15000; Q = '$/caja 15 kilos'; R = 'Región de O''Higgins'; S = 1000 },
  @{ Row = 565; D = 44859; L = 'Segunda'; M = 200; N = 13000; O = 13000; P = 13000; Q = '$/caja 15 kilos'; R = 'Región de O''Higgins'; S = 867 },
  @{ Row = 566; D = 44272; L = 'Primera'; M = 120; N = 16000; O = 17000; P = 16500; Q = '$/caja 15 kilos'; R = 'Región de O''Higgins'; S = 1100 },
  @{ Row = 567; D = 44725; L = 'Especial'; M = 100; N = 20000; O = 20000; P = 20000; Q = '$/caja 15 kilos'; R = 'Región de O''Higgins'; S = 1333 },
  @{ Row = 568; D = 44725; L = 'Primera'; M = 100; N = 17000; O = 17000; P = 17000; Q = '$/caja 15 kilos'; R = 'Región de O''Higgins'; S = 1133 },
  @{ Row = 569; D = 44725; L = 'Segunda'; M = 100; N = 14000; O = 14000; P = 14000; Q = '$/caja 15 kilos'; R = 'Región de O''Higgins'; S = 933 },
  @{ Row = 570; D = 44449; L = 'Especial'; M = 300; N = 20000; O = 20000; P = 20000; Q = '$/caja 15 kilos'; R = 'Provincia de Curicó'; S = 1333 },
  @{ Row = 571; D = 44449; L = 'Primera'; M = 150; N = 14000; O = 14000; P = 14000; Q = '$/caja 15 kilos'; R = 'Provincia de Curicó'; S = 933 },
  @{ Row = 572; D = 44481; L = 'Especial'; M = 200; N = 21000; O = 21000; P = 21000; Q = '$/caja 15 kilos'; R = 'Provincia de Curicó'; S = 1400 },
  @{ Row = 573; D = 44481; L = 'Primera'; M = 500; N = 15000; O = 16000; P = 15500; Q = '$/caja 15 kilos'; R = 'Provincia de Curicó'; S = 1033 },
  @{ Row = 574; D = 44462; L = 'Especial'; M = 200; N = 20000; O = 20000; P = 20000; Q = '$/caja 15 kilos granel'; R = 'Provincia de Curicó'; S = 1333 },
  @{ Row = 575; D = 44991; L = 'Primera'; M = 400; N = 17000; O = 18000; P = 17500; Q = '$/caja 15 kilos'; R = 'Región de O''Higgins'; S = 1167 },
  @{ Row = 576; D = 45076; L = 'Especial'; M = 300; N = 18000; O = 18000; P = 18000; Q = '$/caja 15 kilos'; R = 'Región de O''Higgins'; S = 1200 },
  @{ Row = 577; D = 45076; L = 'Primera'; M = 300; N = 16000; O = 16000; P = 16000; Q = '$/caja 15 kilos'; R = 'Región de O''Higgins'; S = 1067 },
  @{ Row = 578; D = 45076; L = 'Segunda'; M = 300; N = 15000; O = 15000; P = 15000; Q = '$/caja 15 kilos'; R = 'Región de O''Higgins'; S = 1000 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 4).Value = $r.D        # D: Fecha
    $ws.Cells.Item($row, 12).Value = $r.L       # L: Calidad
    $ws.Cells.Item($row, 13).Value = $r.M       # M: Volumen
    $ws.Cells.Item($row, 14).Value = $r.N       # N: Precio minimo
    $ws.Cells.Item($row, 15).Value = $r.O       # O: Precio maximo
    $ws.Cells.Item($row, 16).Value = $r.P       # P: Precio promedio ponderado
    $ws.Cells.Item($row, 17).Value = $r.Q       # Q: Unidad de comercializacion
    $ws.Cells.Item($row, 18).Value = $r.R       # R: Origen
    $ws.Cells.Item($row, 19).Value = $r.S       # S: Precio $/Kg
}